$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1011.4286
$ws.Range("I43").Value = 916
$ws.Range("J43").Value = 1250
$ws.Range("K43").Value = 916
$ws.Range("L43").Value = 1250
$ws.Range("M43").Value = -847
$ws.Range("N43").Value = -1388
$ws.Range("H109").Value = 109960.25
$ws.Range("J109").Value = 109960.25
$ws.Range("L109").Value = 109960.25
$ws.Range("N109").Value = -112734.25
$ws.Range("H116").Value = 20689.785
$ws.Range("I116").Value = 21673.154
$ws.Range("J116").Value = 7906
$ws.Range("K116").Value = 21673.154
$ws.Range("L116").Value = 7906
$ws.Range("M116").Value = -18231.154
$ws.Range("N116").Value = -14790
$ws.Range("H132").Value = 4799.674
$ws.Range("I132").Value = 5532.1816
$ws.Range("K132").Value = 16596.5448
$ws.Range("M132").Value = -14066.5448

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 13498.571
$ws.Range("I3").Value = 3750
$ws.Range("J3").Value = 17398
$ws.Range("K3").Value = 3750
$ws.Range("L3").Value = 17398
$ws.Range("M3").Value = -3635
$ws.Range("N3").Value = -17628
$ws.Range("H38").Value = 4666.3335
$ws.Range("I38").Value = 4666.3335
$ws.Range("K38").Value = 4666.3335
$ws.Range("M38").Value = -4199.3335
$ws.Range("H74").Value = 2463.92
$ws.Range("I74").Value = 2411.389
$ws.Range("K74").Value = 2411.389
$ws.Range("M74").Value = -1537.389
$ws.Range("H77").Value = 2463.92
$ws.Range("I77").Value = 2411.389
$ws.Range("K77").Value = 12056.945
$ws.Range("M77").Value = -7688.945
$ws.Range("H97").Value = 1137
$ws.Range("I97").Value = 755.26666
$ws.Range("K97").Value = 755.26666
$ws.Range("M97").Value = -259.26666
$ws.Range("H122").Value = 6974.6
$ws.Range("I122").Value = 6106.5713
$ws.Range("K122").Value = 18319.7139
$ws.Range("M122").Value = -15869.7139
$ws.Range("H124").Value = 21749.25
$ws.Range("J124").Value = 21749.25
$ws.Range("L124").Value = 21749.25
$ws.Range("N124").Value = -31569.25
$ws.Range("H125").Value = 81145.336
$ws.Range("J125").Value = 81145.336
$ws.Range("L125").Value = 81145.336
$ws.Range("N125").Value = -90985.336

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3785
$ws.Range("I58").Value = 2637.3635
$ws.Range("K58").Value = 2637.3635
$ws.Range("M58").Value = -2434.3635
$ws.Range("H99").Value = 11281.5
$ws.Range("J99").Value = 20566.666
$ws.Range("L99").Value = 20566.666
$ws.Range("N99").Value = -23562.666
$ws.Range("H124").Value = 25970.375
$ws.Range("J124").Value = 23680.428
$ws.Range("L124").Value = 23680.428
$ws.Range("N124").Value = -28590.428
$ws.Range("H126").Value = 11281.5
$ws.Range("J126").Value = 20566.666
$ws.Range("L126").Value = 61699.99800000001
$ws.Range("N126").Value = -66639.99800000001
$ws.Range("H134").Value = 2341.0789
$ws.Range("I134").Value = 1953.4242
$ws.Range("K134").Value = 5860.2726
$ws.Range("M134").Value = -3325.2726
$ws.Range("H136").Value = 3785
$ws.Range("I136").Value = 2637.3635
$ws.Range("K136").Value = 7912.0905
$ws.Range("M136").Value = -5362.0905
$ws.Range("H141").Value = 697632.7
$ws.Range("J141").Value = 760703
$ws.Range("L141").Value = 760703
$ws.Range("N141").Value = -771063

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 34.333332
$ws.Range("I12").Value = 18.333334
$ws.Range("J12").Value = 50.333332
$ws.Range("K12").Value = 55.000002
$ws.Range("L12").Value = 150.999996
$ws.Range("M12").Value = 117.999998
$ws.Range("N12").Value = -496.999996
$ws.Range("H68").Value = 1348.125
$ws.Range("J68").Value = 1363.3334
$ws.Range("L68").Value = 4090.0002
$ws.Range("N68").Value = -5712.0002
$ws.Range("H71").Value = 1348.125
$ws.Range("J71").Value = 1363.3334
$ws.Range("L71").Value = 12270.0006
$ws.Range("N71").Value = -20382.0006

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1046.85
$ws.Range("I97").Value = 729.8333
$ws.Range("J97").Value = 3900
$ws.Range("K97").Value = 729.8333
$ws.Range("L97").Value = 3900
$ws.Range("M97").Value = -233.8333
$ws.Range("N97").Value = -4892
$ws.Range("H105").Value = 68332.664
$ws.Range("J105").Value = 68332.664
$ws.Range("L105").Value = 68332.664
$ws.Range("N105").Value = -75320.664
$ws.Range("H126").Value = 1100.6
$ws.Range("I126").Value = 1085
$ws.Range("K126").Value = 3255
$ws.Range("M126").Value = -785
$ws.Range("H132").Value = 1808.1111
$ws.Range("I132").Value = 1453.4
$ws.Range("J132").Value = 2614.2727
$ws.Range("K132").Value = 4360.200000000001
$ws.Range("L132").Value = 7842.8181
$ws.Range("M132").Value = -1830.200000000001
$ws.Range("N132").Value = -12902.8181

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2606.2
$ws.Range("I7").Value = 2342.2856
$ws.Range("J7").Value = 2837.125
$ws.Range("K7").Value = 2342.2856
$ws.Range("L7").Value = 2837.125
$ws.Range("M7").Value = -2230.2856
$ws.Range("N7").Value = -3061.125
$ws.Range("H22").Value = 2584.12
$ws.Range("J22").Value = 2847.111
$ws.Range("L22").Value = 2847.111
$ws.Range("N22").Value = -3437.111
$ws.Range("H27").Value = 2584.12
$ws.Range("J27").Value = 2847.111
$ws.Range("L27").Value = 2847.111
$ws.Range("N27").Value = -3061.111
$ws.Range("H40").Value = 5774.4287
$ws.Range("I40").Value = 5774.4287
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 5774.4287
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -5638.4287
$ws.Range("N40").ClearContents()
$ws.Range("H55").Value = 494.51852
$ws.Range("I55").Value = 403.52942
$ws.Range("J55").Value = 649.2
$ws.Range("K55").Value = 403.52942
$ws.Range("L55").Value = 649.2
$ws.Range("M55").Value = -230.52942
$ws.Range("N55").Value = -995.2
$ws.Range("H103").Value = 4750
$ws.Range("J103").Value = 4750
$ws.Range("L103").Value = 4750
$ws.Range("N103").Value = -7094
$ws.Range("H126").Value = 2606.2
$ws.Range("I126").Value = 2342.2856
$ws.Range("J126").Value = 2837.125
$ws.Range("K126").Value = 7026.8568
$ws.Range("L126").Value = 8511.375
$ws.Range("M126").Value = -4556.8568
$ws.Range("N126").Value = -13451.375
$ws.Range("H128").Value = 89138.57000000001
$ws.Range("J128").Value = 89138.57000000001
$ws.Range("L128").Value = 89138.57000000001
$ws.Range("N128").Value = -99098.57000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2057.7368
$ws.Range("I122").Value = 2135.647
$ws.Range("K122").Value = 6406.941
$ws.Range("M122").Value = -3956.941
$ws.Range("H136").Value = 5436.5293
$ws.Range("I136").Value = 1955.7307
$ws.Range("K136").Value = 5867.1921
$ws.Range("M136").Value = -3317.1921

Write-Host "Applied all profit-table updates."
